$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.62502065026562
$ws.Range("D2").Value = 9.305633960376841
$ws.Range("E2").Value = 13.85442038235093
$ws.Range("F2").Value = 32.62761883647302
$ws.Range("G2").Value = 33.95986795101835
$ws.Range("H2").Value = 15.6578877538339
$ws.Range("I2").Value = 25.5186437195432
$ws.Range("J2").Value = 10.16688129325169
$ws.Range("K2").Value = 15.18270522333529
$ws.Range("L2").Value = 10.36177555138434
$ws.Range("N2").Value = 17.44698248510248
$ws.Range("O2").Value = 24.5159705517488
$ws.Range("C3").Value = 10.57915925312134
$ws.Range("D3").Value = 9.275413013108222
$ws.Range("E3").Value = 13.84926298767278
$ws.Range("F3").Value = 32.69158447621862
$ws.Range("G3").Value = 34.03836411065929
$ws.Range("H3").Value = 15.70686311111383
$ws.Range("I3").Value = 25.58106363765692
$ws.Range("J3").Value = 10.18425137982849
$ws.Range("K3").Value = 14.76377101289327
$ws.Range("L3").Value = 10.37605248341439
$ws.Range("N3").Value = 17.46131791922075
$ws.Range("O3").Value = 24.59400464789739
$ws.Range("C4").Value = 10.55263694564043
$ws.Range("D4").Value = 9.258049686710134
$ws.Range("E4").Value = 13.84805949154801
$ws.Range("F4").Value = 32.73821529397626
$ws.Range("G4").Value = 34.09654784167172
$ws.Range("H4").Value = 15.73939202794329
$ws.Range("I4").Value = 25.62459629856744
$ws.Range("J4").Value = 10.1959033835155
$ws.Range("K4").Value = 14.50123010104617
$ws.Range("L4").Value = 10.3858429523676
$ws.Range("N4").Value = 17.4718226873251
$ws.Range("O4").Value = 24.64696474030588
$ws.Range("C5").Value = 10.54224871735329
$ws.Range("D5").Value = 9.251278727671275
$ws.Range("E5").Value = 13.84806461318878
$ws.Range("F5").Value = 32.75906429702994
$ws.Range("G5").Value = 34.122760728662
$ws.Range("H5").Value = 15.75326579293449
$ws.Range("I5").Value = 25.64364352243189
$ws.Range("J5").Value = 10.200900184701
$ws.Range("K5").Value = 14.39305264275839
$ws.Range("L5").Value = 10.39009067949049
$ws.Range("N5").Value = 17.47653273848761
$ws.Range("O5").Value = 24.66981326884054
$ws.Range("C6").Value = 10.54054934793037
$ws.Range("D6").Value = 9.250172957859174
$ws.Range("E6").Value = 13.84809544582747
$ws.Range("F6").Value = 32.76263768718828
$ws.Range("G6").Value = 34.12726423382065
$ws.Range("H6").Value = 15.75560684242725
$ws.Range("I6").Value = 25.64688519783872
$ws.Range("J6").Value = 10.2017449211193
$ws.Range("K6").Value = 14.37502245244194
$ws.Range("L6").Value = 10.39081160860017
$ws.Range("N6").Value = 17.47734079758123
$ws.Range("O6").Value = 24.6736836985896
$ws.Range("C7").Value = 10.55249513595472
$ws.Range("D7").Value = 9.257957131062149
$ws.Range("E7").Value = 13.84805755175565
$ws.Range("F7").Value = 32.7384889986621
$ws.Range("G7").Value = 34.09689123732193
$ws.Range("H7").Value = 15.73957663237389
$ws.Range("I7").Value = 25.62484788552706
$ws.Range("J7").Value = 10.19596976536497
$ws.Range("K7").Value = 14.49977579962331
$ws.Range("L7").Value = 10.38589919335915
$ws.Range("N7").Value = 17.47188446917964
$ws.Range("O7").Value = 24.64726775707381
$ws.Range("C8").Value = 10.60887301836678
$ws.Range("D8").Value = 9.294969291993555
$ws.Range("E8").Value = 13.85223585263953
$ws.Range("F8").Value = 32.6481458096413
$ws.Range("G8").Value = 33.98485570391127
$ws.Range("H8").Value = 15.67426428530577
$ws.Range("I8").Value = 25.53908437876747
$ws.Range("J8").Value = 10.17266595168948
$ws.Range("K8").Value = 15.03944428623421
$ws.Range("L8").Value = 10.36648589804236
$ws.Range("N8").Value = 17.45157260365475
$ws.Range("O8").Value = 24.54182793541363
$ws.Range("C9").Value = 10.73203382386035
$ws.Range("D9").Value = 9.376786216476042
$ws.Range("E9").Value = 13.87591623013927
$ws.Range("F9").Value = 32.52946667966158
$ws.Range("G9").Value = 33.84474468547005
$ws.Range("H9").Value = 15.56569227761215
$ws.Range("I9").Value = 25.4122938844315
$ws.Range("J9").Value = 10.13477851336495
$ws.Range("K9").Value = 16.04964591383546
$ws.Range("L9").Value = 10.33652479499248
$ws.Range("N9").Value = 17.42520337597297
$ws.Range("O9").Value = 24.37520085490876
$ws.Range("C10").Value = 10.82966219170772
$ws.Range("D10").Value = 9.442210586909082
$ws.Range("E10").Value = 13.90262487804777
$ws.Range("F10").Value = 32.47806339345487
$ws.Range("G10").Value = 33.79074352203191
$ws.Range("H10").Value = 15.49781973004773
$ws.Range("I10").Value = 25.34447325776463
$ws.Range("J10").Value = 10.11168050283016
$ws.Range("K10").Value = 16.75530672780677
$ws.Range("L10").Value = 10.31942878862601
$ws.Range("N10").Value = 17.4139720362736
$ws.Range("O10").Value = 24.27737650349516
$ws.Range("C11").Value = 10.87550306703491
$ws.Range("D11").Value = 9.473055126639897
$ws.Range("E11").Value = 13.91676525381784
$ws.Range("F11").Value = 32.46246880852904
$ws.Range("G11").Value = 33.77687112867741
$ws.Range("H11").Value = 15.4695266436751
$ws.Range("I11").Value = 25.31913743855441
$ws.Range("J11").Value = 10.10219622003643
$ws.Range("K11").Value = 17.06707418068245
$ws.Range("L11").Value = 10.3127130369815
$ws.Range("N11").Value = 17.410616223148
$ws.Range("O11").Value = 24.23823954046254
$ws.Range("C12").Value = 10.89305635784914
$ws.Range("D12").Value = 9.484884349316069
$ws.Range("E12").Value = 13.92240301435115
$ws.Range("F12").Value = 32.45768410896347
$ws.Range("G12").Value = 33.77315930573329
$ws.Range("H12").Value = 15.45918426043474
$ws.Range("I12").Value = 25.31033760565654
$ws.Range("J12").Value = 10.09875146512776
$ws.Range("K12").Value = 17.18370659940341
$ws.Range("L12").Value = 10.31032205775651
$ws.Range("N12").Value = 17.40959630720878
$ws.Range("O12").Value = 24.2241925851028
$ws.Range("C13").Value = 10.8892674800727
$ws.Range("D13").Value = 9.482330192352888
$ws.Range("E13").Value = 13.92117628232808
$ws.Range("F13").Value = 32.45866473642096
$ws.Range("G13").Value = 33.77389011489114
$ws.Range("H13").Value = 15.46139514723414
$ws.Range("I13").Value = 25.31219746630138
$ws.Range("J13").Value = 10.09948683476609
$ws.Range("K13").Value = 17.15865265781844
$ws.Range("L13").Value = 10.31083023976171
$ws.Range("N13").Value = 17.40980482812903
$ws.Range("O13").Value = 24.22718342148956
$ws.Range("C14").Value = 10.87694335690768
$ws.Range("D14").Value = 9.4740253750431
$ws.Range("E14").Value = 13.9172234204759
$ws.Range("F14").Value = 32.46205270811334
$ws.Range("G14").Value = 33.77653484680202
$ws.Range("H14").Value = 15.46866832093594
$ws.Range("I14").Value = 25.31839754693482
$ws.Range("J14").Value = 10.10190987937698
$ws.Range("K14").Value = 17.0766987312482
$ws.Range("L14").Value = 10.3125132835959
$ws.Range("N14").Value = 17.41052729393706
$ws.Range("O14").Value = 24.23706837836356
$ws.Range("C15").Value = 10.86941944848266
$ws.Range("D15").Value = 9.468957647449736
$ws.Range("E15").Value = 13.91483895128134
$ws.Range("F15").Value = 32.46427388482027
$ws.Range("G15").Value = 33.77835564058314
$ws.Range("H15").Value = 15.47317174935456
$ws.Range("I15").Value = 25.32229874329956
$ws.Range("J15").Value = 10.10341316212008
$ws.Range("K15").Value = 17.02631097392521
$ws.Range("L15").Value = 10.31356399373969
$ws.Range("N15").Value = 17.41100245563743
$ws.Range("O15").Value = 24.24322397320356
$ws.Range("C16").Value = 10.82669419412817
$ws.Range("D16").Value = 9.440216048245107
$ws.Range("E16").Value = 13.9017405625575
$ws.Range("F16").Value = 32.47923933450912
$ws.Range("G16").Value = 33.79186563057409
$ws.Range("H16").Value = 15.49972073625826
$ws.Range("I16").Value = 25.34624010979083
$ws.Range("J16").Value = 10.11232087790258
$ws.Range("K16").Value = 16.73473741514393
$ws.Range("L16").Value = 10.31988899393639
$ws.Range("N16").Value = 17.41422652082957
$ws.Range("O16").Value = 24.28004230833773
$ws.Range("C17").Value = 10.8008418422739
$ws.Range("D17").Value = 9.42285665760496
$ws.Range("E17").Value = 13.89421277815643
$ws.Range("F17").Value = 32.49041571090892
$ws.Range("G17").Value = 33.80289532461418
$ws.Range("H17").Value = 15.51666928491722
$ws.Range("I17").Value = 25.36234102335073
$ws.Range("J17").Value = 10.11804724526591
$ws.Range("K17").Value = 16.55342660673776
$ws.Range("L17").Value = 10.32404066304241
$ws.Range("N17").Value = 17.41665256729112
$ws.Range("O17").Value = 24.30400447008891
$ws.Range("C18").Value = 10.78610746955181
$ws.Range("D18").Value = 9.412974310267105
$ws.Range("E18").Value = 13.89007054856411
$ws.Range("F18").Value = 32.4975772100301
$ws.Range("G18").Value = 33.8102456420202
$ws.Range("H18").Value = 15.52666070283011
$ws.Range("I18").Value = 25.37212099794295
$ws.Range("J18").Value = 10.12143722446702
$ws.Range("K18").Value = 16.44827685186542
$ws.Range("L18").Value = 10.32652851851641
$ws.Range("N18").Value = 17.41821308019359
$ws.Range("O18").Value = 24.31829162577074
$ws.Range("C19").Value = 10.78114222321904
$ws.Range("D19").Value = 9.4096460884242
$ws.Range("E19").Value = 13.8887003649643
$ws.Range("F19").Value = 32.50012786121938
$ws.Range("G19").Value = 33.81290703867538
$ws.Range("H19").Value = 15.53008536388677
$ws.Range("I19").Value = 25.3755214576084
$ws.Range("J19").Value = 10.12260156861013
$ws.Range("K19").Value = 16.41252961233821
$ws.Range("L19").Value = 10.32738804039806
$ws.Range("N19").Value = 17.41876983998851
$ws.Range("O19").Value = 24.32321563853849
$ws.Range("C20").Value = 10.80357995490082
$ws.Range("D20").Value = 9.424694052491263
$ws.Range("E20").Value = 13.89499473460864
$ws.Range("F20").Value = 32.48915008291235
$ws.Range("G20").Value = 33.80161701288846
$ws.Range("H20").Value = 15.51483992532926
$ws.Range("I20").Value = 25.36057331438537
$ws.Range("J20").Value = 10.11742769727495
$ws.Range("K20").Value = 16.57281767415334
$ws.Range("L20").Value = 10.32358837235759
$ws.Range("N20").Value = 17.41637723082431
$ws.Range("O20").Value = 24.30140140162552
$ws.Range("C21").Value = 10.88055806664889
$ws.Range("D21").Value = 9.47646070935304
$ws.Range("E21").Value = 13.91837681424943
$ws.Range("F21").Value = 32.46102716345176
$ws.Range("G21").Value = 33.77571616947122
$ws.Range("H21").Value = 15.46652192668497
$ws.Range("I21").Value = 25.31655486876576
$ws.Range("J21").Value = 10.10119419314002
$ws.Range("K21").Value = 17.10081005451856
$ws.Range("L21").Value = 10.31201480799262
$ws.Range("N21").Value = 17.41030829037533
$ws.Range("O21").Value = 24.23414392237031
$ws.Range("C22").Value = 10.93199512887886
$ws.Range("D22").Value = 9.511158871731851
$ws.Range("E22").Value = 13.9353069303691
$ws.Range("F22").Value = 32.44917894910474
$ws.Range("G22").Value = 33.76777351511093
$ws.Range("H22").Value = 15.43710936414882
$ws.Range("I22").Value = 25.2924163301507
$ws.Range("J22").Value = 10.09143977295907
$ws.Range("K22").Value = 17.43752611262121
$ws.Range("L22").Value = 10.30533731002778
$ws.Range("N22").Value = 17.4078034576796
$ws.Range("O22").Value = 24.19469588039065
$ws.Range("C23").Value = 10.90444280814209
$ws.Range("D23").Value = 9.492562838316978
$ws.Range("E23").Value = 13.92612124709096
$ws.Range("F23").Value = 32.45490487391225
$ws.Range("G23").Value = 33.77118961246373
$ws.Range("H23").Value = 15.45260913249158
$ws.Range("I23").Value = 25.30487558137204
$ws.Range("J23").Value = 10.09656777669309
$ws.Range("K23").Value = 17.25860915474232
$ws.Range("L23").Value = 10.30882026602068
$ws.Range("N23").Value = 17.40900703193578
$ws.Range("O23").Value = 24.21533686173199
$ws.Range("C24").Value = 10.80234165312135
$ws.Range("D24").Value = 9.423863061144505
$ws.Range("E24").Value = 13.89464063376672
$ws.Range("F24").Value = 32.48971998063305
$ws.Range("G24").Value = 33.80219179432972
$ws.Range("H24").Value = 15.51566620795185
$ws.Range("I24").Value = 25.36137086550885
$ws.Range("J24").Value = 10.1177074902118
$ws.Range("K24").Value = 16.56405380491274
$ws.Range("L24").Value = 10.32379253831978
$ws.Range("N24").Value = 17.41650119405474
$ws.Range("O24").Value = 24.30257665732785
$ws.Range("C25").Value = 10.69742046513721
$ws.Range("D25").Value = 9.353695523217027
$ws.Range("E25").Value = 13.86786475969283
$ws.Range("F25").Value = 32.55529519141291
$ws.Range("G25").Value = 33.87408043901905
$ws.Range("H25").Value = 15.59297524443452
$ws.Range("I25").Value = 25.44215173560098
$ws.Range("J25").Value = 10.14419428751203
$ws.Range("K25").Value = 15.78229280932444
$ws.Range("L25").Value = 10.34376471390718
$ws.Range("N25").Value = 17.43090227701698
$ws.Range("O25").Value = 24.41596717794808
